# Edit script: add Sheet4 with combat-log content, clear Sticky Sticks' weakness cell,
# and move the selection on Sheet3.

$wb = $excel.ActiveWorkbook

# --- Sheet3: clear the "Weakness" value for the "Sticky Sticks" row (D17) ---
# (previously "Water"; the sort/weight logic apparently didn't want a weakness there)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("D17").Value = "  "

$values = @(
    "---",
    "Action: Swift Surf",
    "User: FoeGreen",
    "Weight:2",
    "Targets:",
    "FoeGreen",
    "---",
    "Action: Swift Surf",
    "User: FoeGreen",
    "Weight:2",
    "Targets:",
    "FoeGreen",
    "---",
    "Action: Bamboo Bash",
    "User: FoeGreen",
    "Weight:2",
    "Targets:",
    "FoeGreen",
    "---",
    "Action: Bamboo Bash",
    "User: FoeGreen",
    "Weight:2",
    "Targets:",
    "FoeGreen",
    "---",
    "Action: Fire Blitz",
    "User: FoeGreen",
    "Weight:1",
    "Targets:",
    "FoeGreen",
    "FoeRed",
    "FoeBlue",
    "---",
    "Action: Natural Remedy",
    "User: FoeRed",
    "Weight:2",
    "Targets:",
    "FoeRed",
    "---",
    "Action: Natural Remedy",
    "User: FoeRed",
    "Weight:2",
    "Targets:",
    "FoeRed",
    "---",
    "Action: Icicle Blade",
    "User: FoeRed",
    "Weight:3",
    "Targets:",
    "FoeRed",
    "---",
    "Action: Icicle Blade",
    "User: FoeRed",
    "Weight:3",
    "Targets:",
    "FoeRed",
    "---",
    "Action: Icicle Blade",
    "User: FoeRed",
    "Weight:3",
    "Targets:",
    "FoeRed",
    "---",
    "Action: Sticky Sticks",
    "User: FoeRed",
    "Weight:2",
    "Targets:",
    "FoeRed",
    "---",
    "Action: Sticky Sticks",
    "User: FoeRed",
    "Weight:2",
    "Targets:",
    "FoeRed",
    "---",
    "Action: Sticky Sticks",
    "User: FoeRed",
    "Weight:2",
    "Targets:",
    "FoeRed",
    "---",
    "Action: Healing Pulse",
    "User: FoeBlue",
    "Weight:2",
    "Targets:",
    "FoeGreen",
    "FoeRed",
    "FoeBlue",
    "---",
    "Action: Fire Ball",
    "User: FoeBlue",
    "Weight:3",
    "Targets:",
    "FoeBlue",
    "---",
    "Action: Fire Ball",
    "User: FoeBlue",
    "Weight:3",
    "Targets:",
    "FoeBlue",
    "---",
    "Action: Bonfire",
    "User: FoeBlue",
    "Weight:3",
    "Targets:",
    "FoeBlue",
    "---",
    "Action: Bonfire",
    "User: FoeBlue",
    "Weight:3",
    "Targets:",
    "FoeBlue"
)


# --- Add Sheet4 at the end of the tab strip with the combat-log rows ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "Sheet4"

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 2).Value = $values[$i]
}

$ws4.Columns("B:B").ColumnWidth = 19.7
$ws4.Range("D3").Select()
$excel.ActiveWindow.ScrollRow = 103

# --- Restore Sheet3 as the active sheet/selection (unchanged activeTab) ---
$ws3.Activate()
$ws3.Range("D11").Select()

Write-Output "done"
